# Trade #27 (row index 57, 0-based trade id "57") closed at 2026-02-18 00:13:15
# - unknown UNKNOWN +0.000%
# Also a brand-new momentum trade (#86) was opened at 2026-02-18 00:13:09.
# This script reproduces both events plus the summary / strategy-status
# roll-up numbers that shift as a result.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.84    # Current Capital
$summary.Range("B4").Value = 0.94       # Total P&L $
$summary.Range("B5").Value = 0.34       # Total P&L %
$summary.Range("B6").Value = 55         # Total Trades
$summary.Range("B7").Value = 31         # Winning Trades
$summary.Range("B9").Value = 56.36      # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.84   # Capital
$status.Range("D6").Value = 26      # Trades
$status.Range("E6").Value = 0.03    # P&L $
$status.Range("F6").Value = -0.16   # P&L %
$status.Range("G6").Value = 57.69   # Win Rate %

# ---------------------------------------------------------------------------
# 3) All Trades sheet - close out trade #57 (row 58) and append new trade #86
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #57 (MarketMaking) gets closed out early.
$allTrades.Range("G58").Value = 0.76        # Exit Price
$allTrades.Range("H58").Value = "CLOSED"    # Status
$allTrades.Range("I58").Value = 8.5714      # P&L %
$allTrades.Range("J58").Value = 0.06        # P&L $
$allTrades.Range("K58").Value = 99.84       # Capital After
$allTrades.Range("L58").Value = "early_exit" # Exit Reason
$allTrades.Range("M58").Value = 0.13        # Duration (min)

# New trade #86 (momentum, still OPEN) appended as row 87.
$allTrades.Range("A87").Value = 86
$allTrades.Range("B87:C87").NumberFormat = "@"
$allTrades.Range("B87").Value = "2026-02-18"
$allTrades.Range("C87").Value = "00:13:09"
$allTrades.Range("D87").Value = "momentum"
$allTrades.Range("E87").Value = "DOWN"
$allTrades.Range("F87").Value = 0.7
# G87 (Exit Price) stays blank - trade is still open.
$allTrades.Range("H87").Value = "OPEN"
$allTrades.Range("I87").Value = 0
$allTrades.Range("J87").Value = 0
$allTrades.Range("K87").Value = 100
# L87 (Exit Reason) stays blank - trade is still open.
$allTrades.Range("M87").Value = 0
$allTrades.Range("N87").Value = 0
$allTrades.Range("O87").Value = 0
$allTrades.Range("P87").Value = 0.9
$allTrades.Range("Q87").Value = "Downward momentum: -1.980% over 10 samples"

# ---------------------------------------------------------------------------
# 4) momentum strategy sheet - append the same new trade #86 as its row 17
#    (column layout differs: Entry/Exit slippage before Confidence, then
#    Entry Reason / Exit Reason / Duration)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A17").Value = 86
$momentum.Range("B17:C17").NumberFormat = "@"
$momentum.Range("B17").Value = "2026-02-18"
$momentum.Range("C17").Value = "00:13:09"
$momentum.Range("D17").Value = "momentum"
$momentum.Range("E17").Value = "DOWN"
$momentum.Range("F17").Value = 0.7
# G17 (Exit Price) stays blank - trade is still open.
$momentum.Range("H17").Value = "OPEN"
$momentum.Range("I17").Value = 0
$momentum.Range("J17").Value = 0
$momentum.Range("K17").Value = 100
$momentum.Range("L17").Value = 0
$momentum.Range("M17").Value = 0
$momentum.Range("N17").Value = 0.9
$momentum.Range("O17").Value = "Downward momentum: -1.980% over 10 samples"
# P17 (Exit Reason) stays blank - trade is still open.
$momentum.Range("Q17").Value = 0

# ---------------------------------------------------------------------------
# 5) MarketMaking strategy sheet - close out the matching trade #57 (row 29)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G29").Value = 0.76          # Exit Price
$marketMaking.Range("H29").Value = "CLOSED"      # Status
$marketMaking.Range("I29").Value = 8.5714        # P&L %
$marketMaking.Range("J29").Value = 0.06          # P&L $
$marketMaking.Range("K29").Value = 99.84         # Capital After
$marketMaking.Range("P29").Value = "early_exit"  # Exit Reason
$marketMaking.Range("Q29").Value = 0.13          # Duration (min)
